$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.504.83"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.613.89"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.31%  "
$c = $ws.Range("D5")
$c.Value = "'510.80"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "
$c = $ws.Range("D6")
$c.Value = "'154.38"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("E7").Value = "  +0.26%  "
$c = $ws.Range("D8")
$c.Value = "'0.587"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "2.627.83"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "3.070.67"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "60.456.78"
$ws.Range("E15").Value = "  -0.49%  "
$c = $ws.Range("D16")
$c.Value = "'21.63"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "2.625.19"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("E19").Value = "  -0.90%  "
$c = $ws.Range("D20")
$c.Value = "'350.78"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.31%  "
$c = $ws.Range("D21")
$c.Value = "'10.63"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  -0.89%  "
$c = $ws.Range("D23")
$c.Value = "'0.998"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.07%  "
$c = $ws.Range("D24")
$c.Value = "'60.61"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  -1.13%  "
$c = $ws.Range("D27")
$c.Value = "'0.996"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0844"
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  +0.11%  "
$c = $ws.Range("D31")
$c.Value = "'19.46"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D32")
$c.Value = "'1.57"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D33")
$c.Value = "'150.48"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.38%  "
$c = $ws.Range("D34")
$c.Value = "'5.79"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("E36").Value = "  -2.62%  "
$c = $ws.Range("D37")
$c.Value = "'0.880"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D40")
$c.Value = "'36.33"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D41")
$c.Value = "'3.77"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "
$c = $ws.Range("D42")
$c.Value = "'294.22"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.10%  "
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("E44").Value = "  -0.26%  "
$c = $ws.Range("D45")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  -4.53%  "
$c = $ws.Range("D47")
$c.Value = "'19.92"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.49%  "
$c = $ws.Range("D48")
$c.Value = "'4.90"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "2.003.26"
$ws.Range("E51").Value = "  -3.62%  "
